# Apply crypto price/volume updates from the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.712.23'
$ws.Range('E2').Value = '  -1.28%  '

$ws.Range('D3').Value = '2.493.31'
$ws.Range('E3').Value = '  -1.48%  '

$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').Value = '''533.66'
$ws.Range('E5').Value = '  -0.63%  '

$ws.Range('D6').Value = '''136.27'
$ws.Range('E6').Value = '  -3.91%  '

$ws.Range('D7').Value = '''0.996'
$ws.Range('E7').Value = '  +0.14%  '

$ws.Range('D8').Value = '''0.559'
$ws.Range('E8').Value = '  -1.36%  '

$ws.Range('D9').Value = '2.513.67'
$ws.Range('E9').Value = '  -1.06%  '

$ws.Range('D10').Value = '''0.100'
$ws.Range('E10').Value = '  +0.53%  '

$ws.Range('E11').Value = '  -0.31%  '

$ws.Range('E12').Value = '  -1.92%  '

$ws.Range('D13').Value = '''0.346'
$ws.Range('E13').Value = '  -2.96%  '

$ws.Range('D14').Value = '2.929.70'
$ws.Range('E14').Value = '  -1.80%  '

$ws.Range('D15').Value = '''23.08'
$ws.Range('E15').Value = '  -1.60%  '

$ws.Range('D16').Value = '58.684.01'
$ws.Range('E16').Value = '  -1.66%  '

$ws.Range('D17').Value = '''0.0000139'
$ws.Range('E17').Value = '  -1.78%  '

$ws.Range('D18').Value = '2.492.69'
$ws.Range('E18').Value = '  -1.72%  '

$ws.Range('D19').Value = '''10.99'
$ws.Range('E19').Value = '  -0.45%  '

$ws.Range('D20').Value = '''4.23'
$ws.Range('E20').Value = '  -0.77%  '

$ws.Range('D21').Value = '''324.05'
$ws.Range('E21').Value = '  +0.11%  '

$ws.Range('D22').Value = '''0.998'
$ws.Range('E22').Value = '  -0.10%  '

$ws.Range('D23').Value = '''5.84'
$ws.Range('E23').Value = '  -0.17%  '

$ws.Range('D24').Value = '''62.92'
$ws.Range('E24').Value = '  +0.01%  '

$ws.Range('D25').Value = '''0.418'
$ws.Range('E25').Value = '  -1.41%  '

$ws.Range('D26').Value = '''0.165'
$ws.Range('E26').Value = '  -0.29%  '

$ws.Range('D27').Value = '''0.990'
$ws.Range('E27').Value = '  -0.96%  '

$ws.Range('D28').Value = '''7.53'
$ws.Range('E28').Value = '  -4.43%  '

$ws.Range('D29').Value = '''6.79'
$ws.Range('E29').Value = '  -1.58%  '

$ws.Range('D30').Value = '0.0₃0768'
$ws.Range('E30').Value = '  -1.09%  '

$ws.Range('D31').Value = '''1.77'
$ws.Range('E31').Value = '  -1.59%  '

$ws.Range('D32').Value = '''166.04'
$ws.Range('E32').Value = '  +1.04%  '

$ws.Range('E33').Value = '  -0.10%  '

$ws.Range('E34').Value = '  -1.27%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '''1.38'
$ws.Range('E35').Value = '  -4.38%  '

$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').Value = '''18.41'
$ws.Range('E36').Value = '  -0.79%  '

$ws.Range('D37').Value = '''4.09'
$ws.Range('E37').Value = '  -5.10%  '

$ws.Range('D38').Value = '''1.55'
$ws.Range('E38').Value = '  -2.88%  '

$ws.Range('D39').Value = '''36.61'
$ws.Range('E39').Value = '  -0.96%  '

$ws.Range('D40').Value = '''0.821'
$ws.Range('E40').Value = '  +0.77%  '

$ws.Range('D41').Value = '''3.60'
$ws.Range('E41').Value = '  -1.84%  '

$ws.Range('D42').Value = '''5.23'
$ws.Range('E42').Value = '  -3.81%  '

$ws.Range('D43').Value = '''277.70'
$ws.Range('E43').Value = '  -4.79%  '

$ws.Range('D44').Value = '''0.995'
$ws.Range('E44').Value = '  -0.20%  '

$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = '''0.601'
$ws.Range('E45').Value = '  -0.16%  '

$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').Value = '''10.85'
$ws.Range('E46').Value = '  -0.06%  '

$ws.Range('D47').Value = '''125.25'
$ws.Range('E47').Value = '  -0.14%  '

$ws.Range('D48').Value = '''0.0926'

$ws.Range('D49').Value = '''0.0509'
$ws.Range('E49').Value = '  -0.71%  '

$ws.Range('D50').Value = '''0.0221'
$ws.Range('E50').Value = '  -2.14%  '

$ws.Range('D51').Value = '''17.41'
$ws.Range('E51').Value = '  -1.75%  '
